$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.413.53'
$ws.Range("E2").Value = '  +11.83%  '
$ws.Range("D3").Value = '1.878.62'
$ws.Range("E3").Value = '  +8.13%  '
$ws.Range("D4").Value = "'0.9945"
$ws.Range("E4").Value = '  -0.33%  '
$ws.Range("D5").Value = "'250.60"
$ws.Range("E5").Value = '  +4.24%  '
$ws.Range("D6").Value = "'0.9943"
$ws.Range("E6").Value = '  -0.33%  '
$ws.Range("D7").Value = "'0.4966"
$ws.Range("E7").Value = '  +3.70%  '
$ws.Range("D8").Value = "'45.04"
$ws.Range("E8").Value = '  +9.09%  '
$ws.Range("D9").Value = "'0.2852"
$ws.Range("E9").Value = '  +9.98%  '
$ws.Range("D10").Value = "'0.06550"
$ws.Range("E10").Value = '  +6.53%  '
$ws.Range("D11").Value = '1.865.71'
$ws.Range("E11").Value = '  +7.38%  '
$ws.Range("D12").Value = "'17.01"
$ws.Range("E12").Value = '  +5.63%  '
$ws.Range("D13").Value = "'0.07186"
$ws.Range("E13").Value = '  +3.76%  '
$ws.Range("D14").Value = "'0.6668"
$ws.Range("E14").Value = '  +10.73%  '
$ws.Range("D15").Value = "'86.16"
$ws.Range("E15").Value = '  +12.17%  '
$ws.Range("E16").Value = '  +8.47%  '
$ws.Range("D17").Value = '30.404.82'
$ws.Range("E17").Value = '  +11.97%  '
$ws.Range("D18").Value = "'0.9927"
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("D19").Value = "'0.000007517"
$ws.Range("E19").Value = '  +6.22%  '
$ws.Range("D20").Value = "'12.57"
$ws.Range("E20").Value = '  +10.17%  '
$ws.Range("D21").Value = "'0.9948"
$ws.Range("E21").Value = '  -0.32%  '
$ws.Range("D22").Value = '2.095.99'
$ws.Range("E22").Value = '  +7.59%  '
$ws.Range("D23").Value = "'4.712"
$ws.Range("E23").Value = '  +6.67%  '
$ws.Range("D24").Value = "'5.513"
$ws.Range("E24").Value = '  +8.07%  '
$ws.Range("D25").Value = "'8.989"
$ws.Range("E25").Value = '  +7.20%  '
$ws.Range("D26").Value = "'144.08"
$ws.Range("E26").Value = '  +1.53%  '
$ws.Range("D27").Value = "'136.21"
$ws.Range("E27").Value = '  +27.47%  '
$ws.Range("D28").Value = "'16.77"
$ws.Range("E28").Value = '  +10.01%  '
$ws.Range("D29").Value = "'1.942"
$ws.Range("E29").Value = '  +6.76%  '
$ws.Range("D30").Value = "'1.401"
$ws.Range("E30").Value = '  +1.47%  '
$ws.Range("D31").Value = "'4.249"
$ws.Range("E31").Value = '  +7.73%  '
$ws.Range("D32").Value = "'0.08613"
$ws.Range("E32").Value = '  +8.68%  '
$ws.Range("D33").Value = "'3.900"
$ws.Range("E33").Value = '  +6.35%  '
$ws.Range("D34").Value = "'0.05056"
$ws.Range("E34").Value = '  +6.25%  '
$ws.Range("D35").Value = "'1.135"
$ws.Range("E35").Value = '  +12.19%  '
$ws.Range("D36").Value = "'0.6829"
$ws.Range("E36").Value = '  +10.72%  '
$ws.Range("D37").Value = "'2.689"
$ws.Range("E37").Value = '  +3.83%  '
$ws.Range("D38").Value = "'2.309"
$ws.Range("E38").Value = '  +14.39%  '
$ws.Range("E39").Value = '  +8.74%  '
$ws.Range("D40").Value = "'0.9564"
$ws.Range("E40").Value = '  +3.78%  '
$ws.Range("D41").Value = "'0.01628"
$ws.Range("E41").Value = '  +9.43%  '
$ws.Range("D42").Value = "'6.152"
$ws.Range("E42").Value = '  +8.17%  '
$ws.Range("D43").Value = "'103.39"
$ws.Range("E43").Value = '  +4.62%  '
$ws.Range("D44").Value = "'0.9925"
$ws.Range("E44").Value = '  -0.56%  '
$ws.Range("E45").Value = '  +8.90%  '
$ws.Range("D46").Value = "'7.432"
$ws.Range("E46").Value = '  +8.61%  '
$ws.Range("D47").Value = "'0.1250"
$ws.Range("E47").Value = '  +8.60%  '
$ws.Range("D48").Value = "'0.05643"
$ws.Range("E48").Value = '  +5.44%  '
$ws.Range("D49").Value = "'8.345"
$ws.Range("E49").Value = '  +6.91%  '
$ws.Range("D50").Value = "'32.28"
$ws.Range("E50").Value = '  +7.99%  '
$ws.Range("E51").Value = '  +7.89%  '
